$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.256797671318054
$ws.Range("B1").Value = 2.394113540649414
$ws.Range("C1").Value = 4.2109375
$ws.Range("D1").Value = 2.550778150558472
$ws.Range("E1").Value = 1.355120778083801
